# Swap the full content of row 2 <-> row 3, and row 4 <-> row 5.
# (The underlying data rows were reordered; all populated columns in the
# affected rows need to be exchanged so that each row keeps the full set
# of field values belonging to its record.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: columns Y and AA (Startdatum/Slutdatum) hold identical text values
# ("2023-09-10") in every one of these rows, so they are intentionally left
# out here to avoid Excel auto-converting the inline date text into a date
# serial number when the cell is merely round-tripped through .Value2.
# Columns I, K, AT and AY are empty in all four rows, so they are skipped
# too (swapping two empty cells is a no-op, and leaving them alone avoids
# needless churn on these cells).
$columns = @('A','B','C','D','E','F','G','H','P','Q','R','S','T','U','V','W', `
             'Z','AB','AD','AE','AG','AW','AX')

function Swap-Rows([int]$rowA, [int]$rowB) {
    foreach ($col in $columns) {
        $rangeA = $ws.Range("$col$rowA")
        $rangeB = $ws.Range("$col$rowB")
        $valA = $rangeA.Value2
        $valB = $rangeB.Value2
        $rangeA.Value2 = $valB
        $rangeB.Value2 = $valA
    }
}

Swap-Rows 2 3
Swap-Rows 4 5
